$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.535856000000001
$ws.Range("H2").Value = 25.607568
$ws.Range("I2").Value = 0.36987004643386
$ws.Range("J2").Value = 0.36987004643386
$ws.Range("M2").Value = 20.29529466666667
$ws.Range("N2").Value = 60.885884
$ws.Range("O2").Value = 0.4032332285476398
$ws.Range("P2").Value = 0.4032332285476398
$ws.Range("Q2").Value = 173.2377127522347
$ws.Range("R2").Value = 1559.139414770112
$ws.Range("S2").Value = 0.1491438929665908
$ws.Range("T2").Value = 0.1491438929665908
# Row 3
$ws.Range("G3").Value = 8.535856000000001
$ws.Range("H3").Value = 25.607568
$ws.Range("I3").Value = 0.36987004643386
$ws.Range("J3").Value = 0.36987004643386
$ws.Range("O3").Value = 0.1953894087318433
$ws.Range("P3").Value = 0.1953894087318433
$ws.Range("Q3").Value = 83.94351424517335
$ws.Range("R3").Value = 755.4916282065601
$ws.Range("S3").Value = 0.07226868968033132
$ws.Range("T3").Value = 0.07226868968033132
# Row 4
$ws.Range("G4").Value = 8.535856000000001
$ws.Range("H4").Value = 25.607568
$ws.Range("I4").Value = 0.36987004643386
$ws.Range("J4").Value = 0.36987004643386
$ws.Range("M4").Value = 11.81535133333333
$ws.Range("N4").Value = 35.446054
$ws.Range("O4").Value = 0.2347510761885954
$ws.Range("P4").Value = 0.2347510761885954
$ws.Range("Q4").Value = 100.8541375707413
$ws.Range("R4").Value = 907.6872381366719
$ws.Range("S4").Value = 0.08682739145027436
$ws.Range("T4").Value = 0.08682739145027436
# Row 5
$ws.Range("G5").Value = 8.535856000000001
$ws.Range("H5").Value = 25.607568
$ws.Range("I5").Value = 0.36987004643386
$ws.Range("J5").Value = 0.36987004643386
$ws.Range("M5").Value = 8.386535
$ws.Range("N5").Value = 25.159605
$ws.Range("O5").Value = 0.1666262865319216
$ws.Range("P5").Value = 0.1666262865319216
$ws.Range("Q5").Value = 71.58625509896001
$ws.Range("R5").Value = 644.27629589064
$ws.Range("S5").Value = 0.0616300723366635
$ws.Range("T5").Value = 0.06163007233666349
# Row 6
$ws.Range("I6").Value = 0.2236685002562326
$ws.Range("J6").Value = 0.2236685002562326
$ws.Range("M6").Value = 20.29529466666667
$ws.Range("N6").Value = 60.885884
$ws.Range("O6").Value = 0.4032332285476398
$ws.Range("P6").Value = 0.4032332285476398
$ws.Range("Q6").Value = 104.7606308559004
$ws.Range("R6").Value = 942.845677703104
$ws.Range("S6").Value = 0.09019057148272927
$ws.Range("T6").Value = 0.09019057148272928
# Row 7
$ws.Range("I7").Value = 0.2236685002562326
$ws.Range("J7").Value = 0.2236685002562326
$ws.Range("O7").Value = 0.1953894087318433
$ws.Range("P7").Value = 0.1953894087318433
$ws.Range("S7").Value = 0.04370245601700343
$ws.Range("T7").Value = 0.04370245601700343
# Row 8
$ws.Range("I8").Value = 0.2236685002562326
$ws.Range("J8").Value = 0.2236685002562326
$ws.Range("M8").Value = 11.81535133333333
$ws.Range("N8").Value = 35.446054
$ws.Range("O8").Value = 0.2347510761885954
$ws.Range("P8").Value = 0.2347510761885954
$ws.Range("Q8").Value = 60.98870106562487
$ws.Range("R8").Value = 548.898309590624
$ws.Range("S8").Value = 0.05250642114463972
$ws.Range("T8").Value = 0.05250642114463973
# Row 9
$ws.Range("I9").Value = 0.2236685002562326
$ws.Range("J9").Value = 0.2236685002562326
$ws.Range("M9").Value = 8.386535
$ws.Range("N9").Value = 25.159605
$ws.Range("O9").Value = 0.1666262865319216
$ws.Range("P9").Value = 0.1666262865319216
$ws.Range("Q9").Value = 43.28977291165333
$ws.Range("R9").Value = 389.6079562048799
$ws.Range("S9").Value = 0.03726905161186019
$ws.Range("T9").Value = 0.03726905161186019
# Row 10
$ws.Range("G10").Value = 5.335438
$ws.Range("H10").Value = 16.006314
$ws.Range("I10").Value = 0.2311916579666972
$ws.Range("J10").Value = 0.2311916579666973
$ws.Range("M10").Value = 20.29529466666667
$ws.Range("N10").Value = 60.885884
$ws.Range("O10").Value = 0.4032332285476398
$ws.Range("P10").Value = 0.4032332285476398
$ws.Range("Q10").Value = 108.2842863857307
$ws.Range("R10").Value = 974.5585774715761
$ws.Range("S10").Value = 0.093224158655193
$ws.Range("T10").Value = 0.093224158655193
# Row 11
$ws.Range("G11").Value = 5.335438
$ws.Range("H11").Value = 16.006314
$ws.Range("I11").Value = 0.2311916579666972
$ws.Range("J11").Value = 0.2311916579666973
$ws.Range("O11").Value = 0.1953894087318433
$ws.Range("P11").Value = 0.1953894087318433
$ws.Range("Q11").Value = 52.46988887315334
$ws.Range("R11").Value = 472.22899985838
$ws.Range("S11").Value = 0.04517240135384751
$ws.Range("T11").Value = 0.04517240135384752
# Row 12
$ws.Range("G12").Value = 5.335438
$ws.Range("H12").Value = 16.006314
$ws.Range("I12").Value = 0.2311916579666972
$ws.Range("J12").Value = 0.2311916579666973
$ws.Range("M12").Value = 11.81535133333333
$ws.Range("N12").Value = 35.446054
$ws.Range("O12").Value = 0.2347510761885954
$ws.Range("P12").Value = 0.2347510761885954
$ws.Range("Q12").Value = 63.04007448721732
$ws.Range("R12").Value = 567.3606703849559
$ws.Range("S12").Value = 0.05427249051350783
$ws.Range("T12").Value = 0.05427249051350783
# Row 13
$ws.Range("G13").Value = 5.335438
$ws.Range("H13").Value = 16.006314
$ws.Range("I13").Value = 0.2311916579666972
$ws.Range("J13").Value = 0.2311916579666973
$ws.Range("M13").Value = 8.386535
$ws.Range("N13").Value = 25.159605
$ws.Range("O13").Value = 0.1666262865319216
$ws.Range("P13").Value = 0.1666262865319216
$ws.Range("Q13").Value = 44.74583752733
$ws.Range("R13").Value = 402.71253774597
$ws.Range("S13").Value = 0.03852260744414891
$ws.Range("T13").Value = 0.03852260744414891
# Row 14
$ws.Range("G14").Value = 4.044874
$ws.Range("H14").Value = 12.134622
$ws.Range("I14").Value = 0.1752697953432102
$ws.Range("J14").Value = 0.1752697953432102
$ws.Range("M14").Value = 20.29529466666667
$ws.Range("N14").Value = 60.885884
$ws.Range("O14").Value = 0.4032332285476398
$ws.Range("P14").Value = 0.4032332285476398
$ws.Range("Q14").Value = 82.09190971953866
$ws.Range("R14").Value = 738.8271874758481
$ws.Range("S14").Value = 0.07067460544312672
$ws.Range("T14").Value = 0.07067460544312672
# Row 15
$ws.Range("G15").Value = 4.044874
$ws.Range("H15").Value = 12.134622
$ws.Range("I15").Value = 0.1752697953432102
$ws.Range("J15").Value = 0.1752697953432102
$ws.Range("O15").Value = 0.1953894087318433
$ws.Range("P15").Value = 0.1953894087318433
$ws.Range("Q15").Value = 39.77819427119334
$ws.Range("R15").Value = 358.00374844074
$ws.Range("S15").Value = 0.03424586168066102
$ws.Range("T15").Value = 0.03424586168066102
# Row 16
$ws.Range("G16").Value = 4.044874
$ws.Range("H16").Value = 12.134622
$ws.Range("I16").Value = 0.1752697953432102
$ws.Range("J16").Value = 0.1752697953432102
$ws.Range("M16").Value = 11.81535133333333
$ws.Range("N16").Value = 35.446054
$ws.Range("O16").Value = 0.2347510761885954
$ws.Range("P16").Value = 0.2347510761885954
$ws.Range("Q16").Value = 47.79160740906533
$ws.Range("R16").Value = 430.124466681588
$ws.Range("S16").Value = 0.04114477308017345
$ws.Range("T16").Value = 0.04114477308017345
# Row 17
$ws.Range("G17").Value = 4.044874
$ws.Range("H17").Value = 12.134622
$ws.Range("I17").Value = 0.1752697953432102
$ws.Range("J17").Value = 0.1752697953432102
$ws.Range("M17").Value = 8.386535
$ws.Range("N17").Value = 25.159605
$ws.Range("O17").Value = 0.1666262865319216
$ws.Range("P17").Value = 0.1666262865319216
$ws.Range("Q17").Value = 33.92247737159
$ws.Range("R17").Value = 305.30229634431
$ws.Range("S17").Value = 0.029204555139249
$ws.Range("T17").Value = 0.029204555139249
